# chore: update Sheets via scheduled runner
# Refresh market-price-derived figures (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 33
$ws_ALC.Range("H33").Value = 62706.125
$ws_ALC.Range("I33").Value = 66853.2
$ws_ALC.Range("K33").Value = 66853.2
$ws_ALC.Range("M33").Value = -66624.2

# ALC row 52
$ws_ALC.Range("H52").Value = 26596.342
$ws_ALC.Range("J52").Value = 298.97144
$ws_ALC.Range("L52").Value = 896.9143199999999
$ws_ALC.Range("N52").Value = -1216.91432

# ALC row 70
$ws_ALC.Range("H70").Value = 551846.4399999999
$ws_ALC.Range("J70").Value = 112848.89
$ws_ALC.Range("L70").Value = 338546.67
$ws_ALC.Range("N70").Value = -339086.67

# ALC row 73
$ws_ALC.Range("H73").Value = 551846.4399999999
$ws_ALC.Range("J73").Value = 112848.89
$ws_ALC.Range("L73").Value = 338546.67
$ws_ALC.Range("N73").Value = -340418.67

# ALC row 74
$ws_ALC.Range("H74").Value = 20599.125
$ws_ALC.Range("I74").Value = 19133
$ws_ALC.Range("K74").Value = 19133
$ws_ALC.Range("M74").Value = -18197

# ALC row 77
$ws_ALC.Range("H77").Value = 20599.125
$ws_ALC.Range("I77").Value = 19133
$ws_ALC.Range("K77").Value = 95665
$ws_ALC.Range("M77").Value = -90985

# ALC row 112
$ws_ALC.Range("H112").Value = 1756.0625
$ws_ALC.Range("J112").Value = 1756.0625
$ws_ALC.Range("L112").Value = 5268.1875
$ws_ALC.Range("N112").Value = -7484.1875

# ALC row 137
$ws_ALC.Range("H137").Value = 2397.4822
$ws_ALC.Range("I137").Value = 1707.9131
$ws_ALC.Range("J137").Value = 5569.5
$ws_ALC.Range("K137").Value = 5123.7393
$ws_ALC.Range("L137").Value = 16708.5
$ws_ALC.Range("M137").Value = -2573.7393
$ws_ALC.Range("N137").Value = -21808.5

# ALC row 138
$ws_ALC.Range("H138").Value = 6112.0864
$ws_ALC.Range("I138").Value = 3046.0454
$ws_ALC.Range("J138").Value = 7255.356
$ws_ALC.Range("K138").Value = 9138.136200000001
$ws_ALC.Range("L138").Value = 21766.068
$ws_ALC.Range("M138").Value = -3998.136200000001
$ws_ALC.Range("N138").Value = -32046.068

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws_ARM.Range("H32").Value = 8474.947
$ws_ARM.Range("I32").Value = 7328.758
$ws_ARM.Range("J32").Value = 16039.8
$ws_ARM.Range("K32").Value = 7328.758
$ws_ARM.Range("L32").Value = 16039.8
$ws_ARM.Range("M32").Value = -7041.758
$ws_ARM.Range("N32").Value = -16613.8

# ARM row 61
$ws_ARM.Range("H61").Value = 2892.6206
$ws_ARM.Range("I61").Value = 2892.6206
$ws_ARM.Range("K61").Value = 2892.6206
$ws_ARM.Range("M61").Value = -2680.6206

# ARM row 74
$ws_ARM.Range("H74").Value = 2417.6428
$ws_ARM.Range("I74").Value = 2396.074
$ws_ARM.Range("K74").Value = 2396.074
$ws_ARM.Range("M74").Value = -1522.074

# ARM row 77
$ws_ARM.Range("H77").Value = 2417.6428
$ws_ARM.Range("I77").Value = 2396.074
$ws_ARM.Range("K77").Value = 11980.37
$ws_ARM.Range("M77").Value = -7612.370000000001

# ARM row 102
$ws_ARM.Range("H102").Value = 1210.875
$ws_ARM.Range("I102").Value = 1156.2778
$ws_ARM.Range("K102").Value = 1156.2778
$ws_ARM.Range("M102").Value = 465.7221999999999

# ARM row 132
$ws_ARM.Range("H132").Value = 2781.6667
$ws_ARM.Range("I132").Value = 2479.4263
$ws_ARM.Range("K132").Value = 7438.2789
$ws_ARM.Range("M132").Value = -4908.2789

# ARM row 136
$ws_ARM.Range("H136").Value = 2892.6206
$ws_ARM.Range("I136").Value = 2892.6206
$ws_ARM.Range("K136").Value = 8677.861800000001
$ws_ARM.Range("M136").Value = -6127.861800000001

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 22
$ws_BSM.Range("H22").Value = 254.33333
$ws_BSM.Range("I22").Value = 254.33333
$ws_BSM.Range("K22").Value = 254.33333
$ws_BSM.Range("M22").Value = -81.33332999999999

# BSM row 94
$ws_BSM.Range("H94").Value = 1212.8158
$ws_BSM.Range("I94").Value = 1212.8158
$ws_BSM.Range("K94").Value = 1212.8158
$ws_BSM.Range("M94").Value = -761.8158000000001

# BSM row 99
$ws_BSM.Range("H99").Value = 994.5
$ws_BSM.Range("I99").Value = 991.3333
$ws_BSM.Range("J99").Value = 997.6667
$ws_BSM.Range("K99").Value = 991.3333
$ws_BSM.Range("L99").Value = 997.6667
$ws_BSM.Range("M99").Value = 506.6667
$ws_BSM.Range("N99").Value = -3993.6667

# BSM row 134
$ws_BSM.Range("H134").Value = 27896.725
$ws_BSM.Range("I134").Value = 2971
$ws_BSM.Range("K134").Value = 8913
$ws_BSM.Range("M134").Value = -6378

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 21
$ws_CRP.Range("H21").Value = 2100
$ws_CRP.Range("J21").Value = 2400
$ws_CRP.Range("L21").Value = 2400
$ws_CRP.Range("N21").Value = -2870

# CRP row 31
$ws_CRP.Range("H31").Value = 36581.62
$ws_CRP.Range("I31").Value = 1388
$ws_CRP.Range("K31").Value = 1388
$ws_CRP.Range("M31").Value = -1093

# CRP row 34
$ws_CRP.Range("H34").Value = 36581.62
$ws_CRP.Range("I34").Value = 1388
$ws_CRP.Range("K34").Value = 1388
$ws_CRP.Range("M34").Value = -1186

# CRP row 58
$ws_CRP.Range("H58").Value = 1511.5358
$ws_CRP.Range("I58").Value = 1373.9615
$ws_CRP.Range("K58").Value = 1373.9615
$ws_CRP.Range("M58").Value = -1170.9615

# CRP row 99
$ws_CRP.Range("H99").Value = 6492.636
$ws_CRP.Range("I99").Value = 5917
$ws_CRP.Range("K99").Value = 5917
$ws_CRP.Range("M99").Value = -4419

# CRP row 107
$ws_CRP.Range("H107").Value = 1000
$ws_CRP.Range("I107").Value = 0
$ws_CRP.Range("K107").Value = 0
$ws_CRP.Range("M107").ClearContents()

# CRP row 126
$ws_CRP.Range("H126").Value = 6492.636
$ws_CRP.Range("I126").Value = 5917
$ws_CRP.Range("K126").Value = 17751
$ws_CRP.Range("M126").Value = -15281

# CRP row 136
$ws_CRP.Range("H136").Value = 1511.5358
$ws_CRP.Range("I136").Value = 1373.9615
$ws_CRP.Range("K136").Value = 4121.8845
$ws_CRP.Range("M136").Value = -1571.8845

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 56
$ws_CUL.Range("H56").Value = 7333
$ws_CUL.Range("I56").Value = 7333
$ws_CUL.Range("K56").Value = 7333
$ws_CUL.Range("M56").Value = -6803

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 46
$ws_GSM.Range("H46").Value = 45739.8
$ws_GSM.Range("I46").Value = 41000
$ws_GSM.Range("J46").Value = 48899.668
$ws_GSM.Range("K46").Value = 41000
$ws_GSM.Range("L46").Value = 48899.668
$ws_GSM.Range("M46").Value = -40844
$ws_GSM.Range("N46").Value = -49211.668

# GSM row 64
$ws_GSM.Range("H64").Value = 61468.5
$ws_GSM.Range("J64").Value = 61468.5
$ws_GSM.Range("L64").Value = 61468.5
$ws_GSM.Range("N64").Value = -61964.5

# GSM row 67
$ws_GSM.Range("H67").Value = 61468.5
$ws_GSM.Range("J67").Value = 61468.5
$ws_GSM.Range("L67").Value = 61468.5
$ws_GSM.Range("N67").Value = -63184.5

# GSM row 102
$ws_GSM.Range("H102").Value = 3106.6562
$ws_GSM.Range("I102").Value = 1642.25
$ws_GSM.Range("J102").Value = 7499.875
$ws_GSM.Range("K102").Value = 1642.25
$ws_GSM.Range("L102").Value = 7499.875
$ws_GSM.Range("M102").Value = -20.25
$ws_GSM.Range("N102").Value = -10743.875

# GSM row 113
$ws_GSM.Range("H113").Value = 597599.8
$ws_GSM.Range("I113").Value = 1253337.8
$ws_GSM.Range("J113").Value = 14721.667
$ws_GSM.Range("K113").Value = 1253337.8
$ws_GSM.Range("L113").Value = 14721.667
$ws_GSM.Range("M113").Value = -1251167.8
$ws_GSM.Range("N113").Value = -19061.667

# GSM row 122
$ws_GSM.Range("H122").Value = 3107.9707
$ws_GSM.Range("I122").Value = 2780.1667
$ws_GSM.Range("K122").Value = 8340.500100000001
$ws_GSM.Range("M122").Value = -5890.500100000001

# GSM row 126
$ws_GSM.Range("H126").Value = 4436.75
$ws_GSM.Range("I126").Value = 0
$ws_GSM.Range("K126").Value = 0
$ws_GSM.Range("M126").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws_LTW.Range("H22").Value = 1000
$ws_LTW.Range("I22").Value = 1000
$ws_LTW.Range("K22").Value = 1000
$ws_LTW.Range("M22").Value = -705

# LTW row 24
$ws_LTW.Range("H24").Value = 0
$ws_LTW.Range("J24").Value = 0
$ws_LTW.Range("L24").Value = 0
$ws_LTW.Range("N24").ClearContents()

# LTW row 27
$ws_LTW.Range("H27").Value = 1000
$ws_LTW.Range("I27").Value = 1000
$ws_LTW.Range("K27").Value = 1000
$ws_LTW.Range("M27").Value = -893

# LTW row 61
$ws_LTW.Range("H61").Value = 19001.75
$ws_LTW.Range("I61").Value = 19001.75
$ws_LTW.Range("K61").Value = 19001.75
$ws_LTW.Range("M61").Value = -18799.75

# LTW row 99
$ws_LTW.Range("H99").Value = 43418.832
$ws_LTW.Range("I99").Value = 36165.5
$ws_LTW.Range("J99").Value = 57925.5
$ws_LTW.Range("K99").Value = 36165.5
$ws_LTW.Range("L99").Value = 57925.5
$ws_LTW.Range("M99").Value = -33170.5
$ws_LTW.Range("N99").Value = -63915.5

# LTW row 113
$ws_LTW.Range("H113").Value = 19001.75
$ws_LTW.Range("I113").Value = 19001.75
$ws_LTW.Range("K113").Value = 19001.75
$ws_LTW.Range("M113").Value = -16831.75

# LTW row 122
$ws_LTW.Range("H122").Value = 421431.62
$ws_LTW.Range("I122").Value = 717344.2
$ws_LTW.Range("K122").Value = 2152032.6
$ws_LTW.Range("M122").Value = -2149582.6

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 94
$ws_WVR.Range("H94").Value = 45000
$ws_WVR.Range("I94").Value = 45000
$ws_WVR.Range("K94").Value = 45000
$ws_WVR.Range("M94").Value = -44099

# WVR row 107
$ws_WVR.Range("H107").Value = 910.9048
$ws_WVR.Range("I107").Value = 1100.2142
$ws_WVR.Range("K107").Value = 3300.6426
$ws_WVR.Range("M107").Value = -1380.6426

# WVR row 122
$ws_WVR.Range("H122").Value = 37040610
$ws_WVR.Range("I122").Value = 90911130
$ws_WVR.Range("K122").Value = 272733390
$ws_WVR.Range("M122").Value = -272730940

# WVR row 126
$ws_WVR.Range("H126").Value = 1958.5
$ws_WVR.Range("I126").Value = 1954.7273
$ws_WVR.Range("K126").Value = 5864.1819
$ws_WVR.Range("M126").Value = -3394.1819

# WVR row 132
$ws_WVR.Range("H132").Value = 26867.6
$ws_WVR.Range("I132").Value = 1287.2903
$ws_WVR.Range("J132").Value = 114977.555
$ws_WVR.Range("K132").Value = 3861.8709
$ws_WVR.Range("L132").Value = 344932.665
$ws_WVR.Range("M132").Value = -1331.8709
$ws_WVR.Range("N132").Value = -349992.665
